# Apply the commit's text edits to the presentation.
# "VM Watcher" -> "VM Cycler" and "Host Watcher" -> "Host Cycler"
# (renaming recurring background watcher components to "cycler" components).

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq "VM Watcher") {
                    $tr.Text = "VM Cycler"
                } elseif ($tr.Text -eq "Host Watcher") {
                    $tr.Text = "Host Cycler"
                }
            }
        }
    }
}
